$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6544943820224719
$wsSummary.Range("C2").Value = 0.5976331360946746
$wsSummary.Range("D2").Value = 0.9456928838951311
$wsSummary.Range("E2").Value = 0.7324147933284989
$wsSummary.Range("F2").Value = 0.8470311975847031
$wsSummary.Range("G2").Value = 0.9249735822472702
$wsSummary.Range("H2").Value = 0.7804955883796939
$wsSummary.Range("I2").Value = 505
$wsSummary.Range("J2").Value = 340
$wsSummary.Range("K2").Value = 194
$wsSummary.Range("L2").Value = 29

# --- Sheet: Classification Report ---
$wsReport = $wb.Worksheets.Item("Classification Report")
$wsReport.Range("B2").Value = 0.8699551569506726
$wsReport.Range("C2").Value = 0.3632958801498127
$wsReport.Range("D2").Value = 0.512549537648613

$wsReport.Range("B3").Value = 0.5976331360946746
$wsReport.Range("C3").Value = 0.9456928838951311
$wsReport.Range("D3").Value = 0.7324147933284989

$wsReport.Range("B4").Value = 0.6544943820224719
$wsReport.Range("C4").Value = 0.6544943820224719
$wsReport.Range("D4").Value = 0.6544943820224719
$wsReport.Range("E4").Value = 0.6544943820224719

$wsReport.Range("B5").Value = 0.7337941465226736
$wsReport.Range("C5").Value = 0.6544943820224719
$wsReport.Range("D5").Value = 0.6224821654885559

$wsReport.Range("B6").Value = 0.7337941465226737
$wsReport.Range("C6").Value = 0.6544943820224719
$wsReport.Range("D6").Value = 0.6224821654885559

# --- Sheet: Confusion Matrix ---
$wsMatrix = $wb.Worksheets.Item("Confusion Matrix")
$wsMatrix.Range("B2").Value = 194
$wsMatrix.Range("C2").Value = 340
$wsMatrix.Range("B3").Value = 29
$wsMatrix.Range("C3").Value = 505
